# Adds the "index and match" worksheet (with Table2: StudentID / FullName /
# Index / Match, using INDEX()/MATCH() formulas) after the existing
# "excel formula" sheet, and makes it the active tab - matching the
# "Add files via upload" commit.

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)

# --- create the new worksheet, placed after the existing one ---------------
$ws = $wb.Worksheets.Add($null, $firstSheet)
$ws.Name = "index and match"

# --- header row --------------------------------------------------------
$ws.Range("A1").Value = "StudentID"
$ws.Range("B1").Value = "FullName"
$ws.Range("C1").Value = "Index"
$ws.Range("D1").Value = "Match"

# --- data rows -----------------------------------------------------------
$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = "Alice"

$ws.Range("A3").Value = 1006
$ws.Range("B3").Value = "Bob"

$ws.Range("A4").Value = 1012
$ws.Range("B4").Value = "Charlie"

$ws.Range("A5").Value = 1018
$ws.Range("B5").Value = "David"

$ws.Range("A6").Value = 1024
$ws.Range("B6").Value = "Emily"

$ws.Range("A7").Value = 1030
$ws.Range("B7").Value = "Frank"

# --- calculated columns ---------------------------------------------------
$ws.Range("C2:C6").Formula = "=INDEX(B2:B10,1)"
$ws.Range("C7").Formula = "=C4"
$ws.Range("D2:D7").Formula = "=MATCH(B2,B2:B7)"

# --- turn the range into an actual table (ListObject) ---------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:D7"), $null, 1)
$lo.Name = "Table2"
$lo.TableStyle = "TableStyleLight9"

# --- column widths (characters; engine stores width in coarser units) -----
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 9.5

# --- light-gray "medium" borders, mirroring the table-style grid -----------
# NOTE: Range.Borders(xlEdgeTop/Bottom) on a multi-row range only paints the
# OUTER edge of that range (matching real Excel semantics), so every row
# that needs its own top+bottom line is looped one row at a time.
$gray = 13421772  # RGB(204,204,204) == #CCCCCC

function Set-Edge($range, $edgeIndex) {
    $b = $range.Borders.Item($edgeIndex)
    $b.Weight = -4138
    $b.Color = $gray
}

function Set-Box($range, [bool]$left, [bool]$right, [bool]$top, [bool]$bottom) {
    if ($left)   { Set-Edge $range 7 }
    if ($right)  { Set-Edge $range 10 }
    if ($top)    { Set-Edge $range 8 }
    if ($bottom) { Set-Edge $range 9 }
}

# header row (A1 / B1:C1 / D1) - single row, safe to set directly
Set-Box $ws.Range("A1")   $false $true  $false $true
Set-Box $ws.Range("B1:C1") $true  $true  $false $true
Set-Box $ws.Range("D1")   $true  $false $false $true

# data rows 2-6: every row needs its own top+bottom, so loop per row
for ($r = 2; $r -le 6; $r++) {
    Set-Box $ws.Cells.Item($r, 1) $false $true  $true $true   # column A
    Set-Box $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 3)) $true $true $true $true  # columns B:C
    Set-Box $ws.Cells.Item($r, 4) $true  $false $true $true   # column D
}

# last row (A7 / B7:C7 / D7) - single row, safe to set directly
Set-Box $ws.Range("A7")   $false $true  $true $false
Set-Box $ws.Range("B7:C7") $true  $true  $true $false
Set-Box $ws.Range("D7")   $true  $false $true $false

# right-align the StudentID column header + data + last row (matches the
# original table's "Student ID" column alignment)
$ws.Range("A1").HorizontalAlignment = -4152
$ws.Range("A2:A6").HorizontalAlignment = -4152
$ws.Range("A7").HorizontalAlignment = -4152

# --- selection / active tab -----------------------------------------------
$ws.Range("D3").Select()
$ws.Activate()
